$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.226.79'
$ws.Range('E2').Value = '  -2.82%  '
$ws.Range('D3').Value = '2.281.80'
$ws.Range('E3').Value = '  -2.71%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '493.13'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -2.21%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '127.03'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.59%  '
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.528'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -1.82%  '
$ws.Range('D9').Value = '2.282.82'
$ws.Range('E9').Value = '  -3.09%  '
$ws.Range('E10').Value = '  -3.41%  '
$ws.Range('E11').Value = '  +0.28%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.321'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.48%  '
$ws.Range('E13').Value = '  -3.24%  '
$ws.Range('D14').Value = '2.685.58'
$ws.Range('E14').Value = '  -2.77%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '21.52'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.49%  '
$ws.Range('D16').Value = '54.057.16'
$ws.Range('E16').Value = '  -3.02%  '
$ws.Range('E17').Value = '  -2.35%  '
$ws.Range('D18').Value = '2.265.28'
$ws.Range('E18').Value = '  -2.12%  '
$ws.Range('E19').Value = '  -0.40%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '301.40'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -3.09%  '
$ws.Range('E22').Value = '  +2.60%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.998'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.39'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -2.04%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '63.70'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -2.45%  '
$ws.Range('E26').Value = '  +0.34%  '
$ws.Range('E27').Value = '  +0.50%  '
$ws.Range('D28').Value = '2.368.40'
$ws.Range('E28').Value = '  -3.37%  '
$ws.Range('E29').Value = '  +1.79%  '
$ws.Range('E30').Value = '  +0.32%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '165.14'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -3.69%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.59'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.67%  '
$ws.Range('D33').Value = '0.0₃0680'
$ws.Range('E33').Value = '  -3.28%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.87'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +1.91%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('E36').Value = '  +0.22%  '
$ws.Range('E37').Value = '  +0.51%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '17.55'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.75%  '
$ws.Range('E39').Value = '  +1.07%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.866'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +4.50%  '
$ws.Range('E41').Value = '  -0.50%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '35.38'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.81%  '
$ws.Range('E43').Value = '  +0.75%  '
$ws.Range('E44').Value = '  +0.94%  '
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '125.68'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.58%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.78'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.10%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.543'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.18%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '236.37'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.93%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0478'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.72%  '
